$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated activity data (runs, balls, fours, sixes) for each match row
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("C3").Value = 72
$ws.Range("D3").Value = 53
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 2

$ws.Range("C4").Value = 62
$ws.Range("D4").Value = 49
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 1

$ws.Range("C5").Value = 65
$ws.Range("D5").Value = 51
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 3

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
